{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The edit removes the \"Requirements:\" section entirely: the heading\n// paragraph plus its four list-item paragraphs (2 charging modes / 200 mA /\n// 600 mA +- 100 mA / Switching done via pin 13...).\nconst textsToDelete = [\n  \"Requirements:\",\n  \"2 charging modes (from battery perspective)\",\n  \"200 mA\",\n  \"600 mA +- 100 mA\",\n  \"Switching done via pin 13, connected to PB4 on microcontroller\"\n];\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (textsToDelete.includes(para.text)) {\n    toDelete.push(para);\n  }\n}\n\n// Delete from the bottom up so earlier indices/ranges remain valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The edit removes the \"Requirements:\" section entirely: the heading\n# paragraph plus its four list-item paragraphs (2 charging modes / 200 mA /\n# 600 mA +- 100 mA / Switching done via pin 13...).\n$startText = \"Requirements:\"\n$endText = \"Switching done via pin 13, connected to PB4 on microcontroller\"\n\n$startPara = $null\n$endPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($null -eq $startPara -and $t -eq $startText) {\n        $startPara = $p\n    }\n    if ($t -eq $endText) {\n        $endPara = $p\n    }\n    if ($null -ne $startPara -and $null -ne $endPara) {\n        break\n    }\n}\n\nif ($null -ne $startPara -and $null -ne $endPara) {\n    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $range.Delete()\n}\n"}
